$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 12:12:39"
$wsZhCn.Range("H3").Value = "2016-03-22 12:13:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 12:12:47"
$wsDeDe.Range("H3").Value = "2016-03-22 12:13:40"
